# Aggiorno file need_to_buy.xlsx da R
# The source feed (column R) advanced by one day: every existing row's
# data shifted up by one position (row N now holds what used to be in
# row N+1), the date in column A was incremented accordingly, and the
# derived/refreshed columns (fcs, buy_BEE_MWH, sell_lago_MWH and
# need_to_buy_MW) were recalculated for the new day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 2-15 (columns A:F) after the refresh.
$data = @(
    @{ Row = 2;  A = 45897; B = 5341.91604070457; C = 4370.93895847995; D = 4704; E = 5689.095951; F = 0.588286198974136 }
    @{ Row = 3;  A = 45898; B = 5341.91604070659; C = 3701.08477377648; D = 1944; E = 5689.095951; F = 87.6776951695785 }
    @{ Row = 4;  A = 45899; B = 975.205470885634; C = 1626.40385460727; D = 1944; E = 1881.064742; F = 24.5109635717347 }
    @{ Row = 5;  A = 45900; B = 862.92289373491;  C = 1576.56293257249; D = 1944; E = 1758.693353; F = 22.0138913265657 }
    @{ Row = 6;  A = 45901; B = 5594.48067261627; C = 4921.16010241;    D = 2952; E = 6203.586309; F = 107.427739116405 }
    @{ Row = 7;  A = 45902; B = 5601.17759246688; C = 5122.20544265882; D = 2952; E = 6398.767756; F = 123.658150257997 }
    @{ Row = 8;  A = 45903; B = 5784.11073458408; C = 5210.22910014964; D = 2952; E = 6398.767756; F = 119.703588398565 }
    @{ Row = 9;  A = 45904; B = 5812.79253171584; C = 5201.0390995687;  D = 2952; E = 6445.459654; F = 120.071092577203 }
    @{ Row = 10; A = 45905; B = 5812.79253171584; C = 4455.72155531159; D = 2952; E = 6445.459654; F = 89.0161948998228 }
    @{ Row = 11; A = 45906; B = 1096.28225605742; C = 2192.07944272347; D = 2952; E = 2113.120328; F = 10.7048964444188 }
    @{ Row = 12; A = 45907; B = 970.493609902267; C = 2140.07235928172; D = 2952; E = 1974.379807; F = 7.99827318247742 }
    @{ Row = 13; A = 45908; B = 6110.60514061168; C = 5038.10823552584; D = 2952; E = 6751.177035; F = 113.611672079757 }
    @{ Row = 14; A = 45909; B = 6110.60514061168; C = 5232.87770721291; D = 2952; E = 6751.177035; F = 121.727066733385 }
    @{ Row = 15; A = 45910; B = 6110.60514061168; C = 5229.48947399595; D = 2952; E = 6751.177035; F = 121.585890349345 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Range("A$r").Value = $entry.A
    $ws.Range("B$r").Value = $entry.B
    $ws.Range("C$r").Value = $entry.C
    $ws.Range("D$r").Value = $entry.D
    $ws.Range("E$r").Value = $entry.E
    $ws.Range("F$r").Value = $entry.F
}
